$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Test")
$ws.Activate()

# Duplicate the formatting of row 5 into row 6 for the new test case LOG-03
$ws.Range("B5:L5").Copy()
$ws.Range("B6:L6").PasteSpecial(-4122)
$ws.Rows.Item(6).RowHeight = 63

# Fill in the new LOG-03 / LOG-TC-03 test case data
$ws.Range("B6").Value = "LOG-03"
$ws.Range("C6").Value = "Verify error when email not registered"
$ws.Range("D6").Value = "LOG-TC-03"
$ws.Range("E6").Value = "Login Module"
$ws.Range("F6").Value = "LOG-TS-03"
$ws.Range("G6").Value = "Login with unregistered email"
$ws.Range("H6").Value = "On login page "
$ws.Range("I6").Value = "1. Enter unregistered email `n2. Enter any password `n3. Click Login"
$ws.Range("J6").Value = "unkown@gmail.com`nPassword123"
$ws.Range("K6").Value = "Error message displayed"
$ws.Range("L6").Value = "High"

# Add the mailto hyperlink on the test data cell, matching the style of the existing rows
$ws.Hyperlinks.Add($ws.Range("J6"), "mailto:unkown@gmail.com`nPassword123")

# Update selection / active cell to reflect the new entry
$ws.Range("J6").Select()
